$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Duel Decks Ajani vs. Nicol Bolas (DDH)',
    'Akki Coalflinger',
    'Allosaurus Rider',
    'Ambush Commander',
    'Boggart Shenanigans',
    'Clickslither',
    'Elvish Eulogist',
    'Elvish Harbinger',
    'Elvish Promenade',
    'Elvish Warrior',
    'Emberwilde Augur',
    'Flamewave Invoker',
    'Forest',
    'Forest',
    'Forest',
    'Forest',
    'Forgotten Cave',
    'Gempalm Incinerator',
    'Gempalm Strider',
    'Giant Growth',
    'Goblin Burrows',
    'Goblin Cohort',
    'Goblin Matron',
    'Goblin Ringleader',
    'Goblin Sledder',
    'Goblin Warchief',
    'Harmonize',
    'Heedless One',
    'Ib Halfheart, Goblin Tactician',
    'Imperious Perfect',
    'Llanowar Elves',
    'Lys Alana Huntmaster',
    'Mogg Fanatic',
    'Mogg War Marshal',
    'Moonglove Extract',
    'Mountain',
    'Mountain',
    'Mountain',
    'Mountain',
    'Mudbutton Torchrunner',
    'Raging Goblin',
    'Reckless One',
    'Siege-Gang Commander',
    'Skirk Drill Sergeant',
    'Skirk Fire Marshal',
    'Skirk Prospector',
    'Skirk Shaman',
    'Slate of Ancestry',
    'Spitting Earth',
    'Stonewood Invoker',
    'Sylvan Messenger',
    'Tarfire',
    'Tar Pitcher',
    'Timberwatch Elf',
    'Tranquil Thicket',
    'Voice of the Woods',
    'Wellwisher',
    'Wildsize',
    'Wirewood Herald',
    'Wirewood Lodge',
    'Wirewood Symbiote',
    'Wood Elves',
    'Wren''s Run Vanquisher',
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Remove the now-extra row 64 (sheet shrank from 64 to 63 rows)
$ws.Rows.Item(64).Delete()

Write-Output "Done"